$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A86").Value = "2025-08"
$ws.Range("B86").Value = 7
$ws.Range("C86").Value = 231
$ws.Range("D86").Value = 3.03030303030303
